$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.953.72"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.985.07"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.33"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.98"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0800"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.01"
$ws.Range("E12").Value = "  +9.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.20"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.847"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "2.274.50"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").Value = "1.982.62"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "36.856.24"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.30"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.10"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.148"
$ws.Range("E26").Value = "  +4.46%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.21"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +19.19%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  +6.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.56"
$ws.Range("E39").Value = "  -7.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0998"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.62"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.32"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.371.90"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.28"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.30"
$ws.Range("E50").Value = "  +5.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("E51").Value = "  +10.95%  "
